$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-10: combinations of sender (ECs/FAPs/sCs) x target (ECs/FAPs/sCs)
# Columns: A=Sending cluster, B=Ligand symbol(Vcan), C=Receptor symbol(Cd44), D=Target cluster
# E..T numeric stats

$rows = @(
    @{A="ECs"; B="Vcan"; C="Cd44"; D="ECs";  E=2; F=0.6666666666666666; G=1.529781; H=4.589343; I=0.01315047351877542; J=0.01315047351877542; K=3; L=1; M=261.380203; N=784.1406089999999; O=0.6968677182772199; P=0.6968677182772199; Q=399.8544683255429; R=3598.690214929886; S=0.009164140475294031; T=0.009164140475294031},
    @{A="ECs"; B="Vcan"; C="Cd44"; D="FAPs"; E=2; F=0.6666666666666666; G=1.529781; H=4.589343; I=0.01315047351877542; J=0.01315047351877542; K=3; L=1; M=31.999428; N=95.998284; O=0.08531391482826334; P=0.08531391482826335; Q=48.95211696526799; R=440.569052687412; S=0.001121918377732139; T=0.001121918377732139},
    @{A="ECs"; B="Vcan"; C="Cd44"; D="sCs";  E=2; F=0.6666666666666666; G=1.529781; H=4.589343; I=0.01315047351877542; J=0.01315047351877542; K=3; L=1; M=81.69901900000001; N=245.097057; O=0.2178183668945166; P=0.2178183668945167; Q=124.981606984839; R=1124.834462863551; S=0.00286441466574925; T=0.00286441466574925},
    @{A="FAPs"; B="Vcan"; C="Cd44"; D="ECs";  E=3; F=1; G=103.676216; H=311.028648; I=0.8912330150752564; J=0.8912330150752563; K=3; L=1; M=261.380203; N=784.1406089999999; O=0.6968677182772199; P=0.6968677182772199; Q=27098.91038435185; R=243890.1934591666; S=0.6210715176688211; T=0.6210715176688211},
    @{A="FAPs"; B="Vcan"; C="Cd44"; D="FAPs"; E=3; F=1; G=103.676216; H=311.028648; I=0.8912330150752564; J=0.8912330150752563; K=3; L=1; M=31.999428; N=95.998284; O=0.08531391482826334; P=0.08531391482826335; Q=3317.579609204447; R=29858.21648284003; S=0.07603457754026677; T=0.07603457754026677},
    @{A="FAPs"; B="Vcan"; C="Cd44"; D="sCs";  E=3; F=1; G=103.676216; H=311.028648; I=0.8912330150752564; J=0.8912330150752563; K=3; L=1; M=81.69901900000001; N=245.097057; O=0.2178183668945166; P=0.2178183668945167; Q=8470.245140832105; R=76232.20626748893; S=0.1941269198661685; T=0.1941269198661685},
    @{A="sCs";  B="Vcan"; C="Cd44"; D="ECs";  E=3; F=1; G=11.12297; H=33.36891; I=0.09561651140596822; J=0.09561651140596822; K=3; L=1; M=261.380203; N=784.1406089999999; O=0.6968677182772199; P=0.6968677182772199; Q=2907.32415656291; R=26165.91740906619; S=0.06663206013310485; T=0.06663206013310485},
    @{A="sCs";  B="Vcan"; C="Cd44"; D="FAPs"; E=3; F=1; G=11.12297; H=33.36891; I=0.09561651140596822; J=0.09561651140596822; K=3; L=1; M=31.999428; N=95.998284; O=0.08531391482826334; P=0.08531391482826335; Q=355.92867766116; R=3203.35809895044; S=0.008157418910264443; T=0.008157418910264443},
    @{A="sCs";  B="Vcan"; C="Cd44"; D="sCs";  E=3; F=1; G=11.12297; H=33.36891; I=0.09561651140596822; J=0.09561651140596822; K=3; L=1; M=81.69901900000001; N=245.097057; O=0.2178183668945166; P=0.2178183668945167; Q=908.7357373664302; R=8178.62163629787; S=0.02082703236259892; T=0.02082703236259892}
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    foreach ($col in $cols) {
        $cellRef = "$col$r"
        $ws.Range($cellRef).Value = $rowData[$col]
    }
}
